$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48; this pushes the existing rows 48..131
# down to 49..132 (matching the diff, which shows every row from 48 onward
# shifting down by one position) and leaves a fresh blank row 48 to fill in.
$ws.Rows(48).Insert()

# Populate the newly inserted row 48 with the new record.
$ws.Range("A48").Value = 10
$ws.Range("B48").Value = "Vega Modelo de Temuco"
$ws.Range("C48").Value = "La Araucanía"
$ws.Range("D48").Value = 44799
$ws.Range("E48").Value = 9
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100107
$ws.Range("H48").Value = "Otros"
$ws.Range("I48").Value = 100107002
$ws.Range("J48").Value = "Chirimoya"
$ws.Range("K48").Value = "Cultivar IV Región"
$ws.Range("L48").Value = "Primera"
$ws.Range("M48").Value = 55
$ws.Range("N48").Value = 4500
$ws.Range("O48").Value = 4500
$ws.Range("P48").Value = 4500
$ws.Range("Q48").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R48").Value = "Provincia del Elquí"
$ws.Range("S48").Value = 4500
$ws.Range("T48").Value = 1
